$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# G2: "Change Start page" steps -- FullScreen value flips from 0 to 1
$ws.Range("G2").Value = $ws.Range("G2").Value.Replace('FullScreen,<FullScreen value="0"/>', 'FullScreen,<FullScreen value="1"/>')

# H13: VT284-0003 validation -- screenshot check replaced by icon-displayed check
$ws.Range("H13").Value = $ws.Range("H13").Value.Replace('validate_Screenshot=VT284_0003', 'validate_isIconDisplayed=batteryview_xpath,true')

# H14: VT284-0006 validation -- screenshot check replaced by icon-displayed check
$ws.Range("H14").Value = $ws.Range("H14").Value.Replace('validate_Screenshot=VT284_0006', 'validate_isIconDisplayed=batteryview_xpath,true')

# G15: VT284-0012 steps -- drop the TakeScreenshot call and the trailing blank line
$ws.Range("G15").Value = $ws.Range("G15").Value.Replace("TakeScreenshot(VT284_0012);`n", "").TrimEnd("`n")

# H15: VT284-0012 validation -- screenshot check replaced by icon-displayed check
$ws.Range("H15").Value = $ws.Range("H15").Value.Replace('validate_Screenshot=VT284_0012', 'validate_isIconDisplayed=batteryview_xpath,true')

# G18: VT284-0022 steps -- drop the trailing blank line only (content unchanged otherwise)
$ws.Range("G18").Value = $ws.Range("G18").Value.TrimEnd("`n")

# G20: VT284-0037 steps -- drop the TakeScreenshot call, keep the trailing blank line
$ws.Range("G20").Value = $ws.Range("G20").Value.Replace("TakeScreenshot(VT284_0037);`n", "")

# H20: VT284-0037 validation -- screenshot check replaced by icon-not-displayed check
$ws.Range("H20").Value = $ws.Range("H20").Value.Replace('validate_Screenshot=VT284_0037', 'validate_isIconDisplayed=batteryview_xpath,false')

# Update the sheet's saved selection to D1
$ws.Range("D1").Select()
